$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.635.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.92%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.025.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.31%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.36%  "

$ws.Range("E7").Value = "  +0.87%  "

$ws.Range("E8").Value = "  +0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.740"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.176"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000340"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.669.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.064.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.60%  "

$ws.Range("E18").Value = "  +0.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.664.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "440.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "94.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.45%  "

$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.04%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.54%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.78%  "

$ws.Range("E28").Value = "  +1.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "700.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.131"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.77%  "

$ws.Range("E32").Value = "  +2.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +13.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "67.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0908"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.84%  "

$ws.Range("E36").Value = "  -1.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "40.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.70%  "

$ws.Range("E38").Value = "  +5.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +18.71%  "

$ws.Range("E40").Value = "  +0.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("E42").Value = "  +1.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.40%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.39%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.47%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.146"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000280"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +18.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0343"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.13%  "
